# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" and "全部类型" worksheets to match the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F4").Value = 4716
        $ws.Range("F6").Value = 419
        $ws.Range("F7").Value = 1427
        $ws.Range("F8").Value = 933
        $ws.Range("F9").Value = 60
        $ws.Range("F10").Value = 1264
        $ws.Range("F12").Value = 942
        $ws.Range("F14").Value = 71
        $ws.Range("F16").Value = 286
    } else {
        $ws.Range("F4").Value = 4716
        $ws.Range("F6").Value = 419
        $ws.Range("F7").Value = 1427
        $ws.Range("F8").Value = 933
        $ws.Range("F9").Value = 60
        $ws.Range("F10").Value = 1264
        $ws.Range("F12").Value = 943
        $ws.Range("F14").Value = 71
        $ws.Range("F16").Value = 286
    }
}
